# Apply the AHP "Caso Extremo" update:
#   - Matriz_Económico: several pairwise-comparison cells change
#     (5 -> 3, 9 -> 7, and their reciprocals 0.2 -> 1/3, 1/9 -> 1/7)
#   - Pesos_Locales_Económico: recomputed local AHP weights for the
#     "Económico" criterion
#   - Resultados / Ranking_Alternativas: recomputed global weights
#     (= average of the five local weight sheets), which also changes
#     the sort order (and therefore the Alternativa names per row) in
#     Ranking_Alternativas

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Matriz_Económico — direct pairwise-comparison edits
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Matriz_Económico")

$third = 0.3333333333333333

# Row 2 (Baron)
$ws.Range("D2").Value = $third
$ws.Range("F2").Value = $third
$ws.Range("G2").Value = $third
$ws.Range("M2").Value = $third
$ws.Range("N2").Value = 7

# Row 3 (Cordillera)
$ws.Range("D3").Value = $third
$ws.Range("F3").Value = $third
$ws.Range("G3").Value = $third
$ws.Range("M3").Value = $third
$ws.Range("N3").Value = 7

# Row 4 (Esperanza)
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 3
$ws.Range("E4").Value = 3
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 3

# Row 5 (Jean y Marie Thierry)
$ws.Range("D5").Value = $third
$ws.Range("F5").Value = $third
$ws.Range("G5").Value = $third
$ws.Range("M5").Value = $third
$ws.Range("N5").Value = 7

# Row 6 (Laguna Verde)
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = 3
$ws.Range("E6").Value = 3
$ws.Range("J6").Value = 3
$ws.Range("K6").Value = 3

# Row 7 (Las Cañas)
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = 3
$ws.Range("E7").Value = 3
$ws.Range("J7").Value = 3
$ws.Range("K7").Value = 3

# Row 10 (Placeres)
$ws.Range("D10").Value = $third
$ws.Range("F10").Value = $third
$ws.Range("G10").Value = $third
$ws.Range("M10").Value = $third
$ws.Range("N10").Value = 7

# Row 11 (Placilla)
$ws.Range("D11").Value = $third
$ws.Range("F11").Value = $third
$ws.Range("G11").Value = $third
$ws.Range("M11").Value = $third

# Row 13 (Puertas Negras)
$ws.Range("B13").Value = 3
$ws.Range("C13").Value = 3
$ws.Range("E13").Value = 3
$ws.Range("J13").Value = 3
$ws.Range("K13").Value = 3

# Row 14 (Quebrada Verde)
$oneSeventh = 0.1428571428571428
$ws.Range("B14").Value = $oneSeventh
$ws.Range("C14").Value = $oneSeventh
$ws.Range("E14").Value = $oneSeventh
$ws.Range("J14").Value = $oneSeventh

# ---------------------------------------------------------------
# 2) Pesos_Locales_Económico — recomputed local weights (B column)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Pesos_Locales_Económico")

$ws.Range("B2").Value = 0.06859393436079969
$ws.Range("B3").Value = 0.06859393436079969
$ws.Range("B4").Value = 0.1451001396860344
$ws.Range("B5").Value = 0.06859393436079965
$ws.Range("B6").Value = 0.1451001396860344
$ws.Range("B7").Value = 0.1451001396860344
$ws.Range("B8").Value = 0.01428744057464852
$ws.Range("B9").Value = 0.01428744057464852
$ws.Range("B10").Value = 0.06859393436079965
$ws.Range("B11").Value = 0.0697866565851307
$ws.Range("B12").Value = 0.008867907233478526
$ws.Range("B13").Value = 0.1451001396860344
$ws.Range("B14").Value = 0.009419377695460395
$ws.Range("B15").Value = 0.01428744057464852
$ws.Range("B16").Value = 0.01428744057464852

# ---------------------------------------------------------------
# 3) Resultados — recomputed global weights, alphabetical order
#    (names unchanged, only the weight column moves)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Resultados")

$ws.Range("B2").Value = 0.05290648748731419
$ws.Range("B3").Value = 0.05575513621927968
$ws.Range("B4").Value = 0.05410508522144566
$ws.Range("B5").Value = 0.08817071047666492
$ws.Range("B6").Value = 0.05776363253761711
$ws.Range("B7").Value = 0.06105912745907267
$ws.Range("B8").Value = 0.08859478509086688
$ws.Range("B9").Value = 0.02694557311450325
$ws.Range("B10").Value = 0.06234583198288431
$ws.Range("B11").Value = 0.1135151614696607
$ws.Range("B12").Value = 0.1161686867745278
$ws.Range("B13").Value = 0.05761308903197804
$ws.Range("B14").Value = 0.06069805844247664
$ws.Range("B15").Value = 0.06126973961633277
$ws.Range("B16").Value = 0.0430888950753756

# ---------------------------------------------------------------
# 4) Ranking_Alternativas — recomputed global weights, re-sorted
#    descending by weight, so Alternativa (A) and Ranking (C) also
#    shift per row; rewrite A/B for every data row.
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Ranking_Alternativas")

$ws.Range("A2").Value = "Plaza Justicia"
$ws.Range("B2").Value = 0.1161686867745278

$ws.Range("A3").Value = "Placilla"
$ws.Range("B3").Value = 0.1135151614696607

$ws.Range("A4").Value = "Marcelo Mena"
$ws.Range("B4").Value = 0.08859478509086688

$ws.Range("A5").Value = "Jean y Marie Thierry"
$ws.Range("B5").Value = 0.08817071047666492

$ws.Range("A6").Value = "Placeres"
$ws.Range("B6").Value = 0.06234583198288431

$ws.Range("A7").Value = "Reina Isabel 2"
$ws.Range("B7").Value = 0.06126973961633277

$ws.Range("A8").Value = "Las Cañas"
$ws.Range("B8").Value = 0.06105912745907267

$ws.Range("A9").Value = "Quebrada Verde"
$ws.Range("B9").Value = 0.06069805844247664

$ws.Range("A10").Value = "Laguna Verde"
$ws.Range("B10").Value = 0.05776363253761711

$ws.Range("A11").Value = "Puertas Negras"
$ws.Range("B11").Value = 0.05761308903197804

$ws.Range("A12").Value = "Cordillera"
$ws.Range("B12").Value = 0.05575513621927968

$ws.Range("A13").Value = "Esperanza"
$ws.Range("B13").Value = 0.05410508522144566

$ws.Range("A14").Value = "Baron"
$ws.Range("B14").Value = 0.05290648748731419

$ws.Range("A15").Value = "Rodelillo"
$ws.Range("B15").Value = 0.0430888950753756

$ws.Range("A16").Value = "Padre Demian Molokai"
$ws.Range("B16").Value = 0.02694557311450325

# C (Ranking) column is unchanged (1..15 top-to-bottom) — left as-is.
